$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 44
$ws1.Range("F3").Value = 26450
$ws1.Range("F7").Value = 175
$ws1.Range("F8").Value = 538
$ws1.Range("F11").Value = 224
$ws1.Range("F16").Value = 385
$ws1.Range("F20").Value = 35

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 191
$ws2.Range("F8").Value = 111
$ws2.Range("F9").Value = 111
$ws2.Range("F10").Value = 435

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5019
$ws3.Range("F3").Value = 218

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 44
$ws4.Range("F3").Value = 5019
$ws4.Range("F4").Value = 218
$ws4.Range("F5").Value = 26450
$ws4.Range("F13").Value = 175
$ws4.Range("F14").Value = 191
$ws4.Range("F15").Value = 191
$ws4.Range("F17").Value = 111
$ws4.Range("F18").Value = 111
$ws4.Range("F19").Value = 435
$ws4.Range("F20").Value = 538
$ws4.Range("F24").Value = 224
$ws4.Range("F32").Value = 385
$ws4.Range("F38").Value = 35
